{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// worksheet table with its new problem, preserving all run/paragraph\n// formatting (font, size, justification, etc.) by editing the text in\n// place via Range.insertText(..., Word.InsertLocation.replace).\nconst replacements = [\n  [\"520\u00d75=\", \"817\u00d78=\"],\n  [\"249\u00d76=\", \"779\u00d79=\"],\n  [\"442\u00d77=\", \"936\u00d78=\"],\n  [\"377\u00d75=\", \"220\u00d76=\"],\n  [\"382\u00d76=\", \"984\u00d76=\"],\n  [\"468\u00d79=\", \"564\u00d72=\"],\n  [\"610\u00d73=\", \"944\u00d72=\"],\n  [\"374\u00d79=\", \"555\u00d75=\"],\n  [\"171\u00d79=\", \"357\u00d76=\"],\n  [\"966\u00d79=\", \"982\u00d74=\"],\n  [\"942\u00d79=\", \"595\u00d76=\"],\n  [\"710\u00d76=\", \"931\u00d75=\"],\n  [\"135\u00d76=\", \"285\u00d77=\"],\n  [\"576\u00d74=\", \"816\u00d73=\"],\n  [\"429\u00d78=\", \"755\u00d75=\"],\n  [\"383\u00d73=\", \"825\u00d74=\"],\n  [\"502\u00d79=\", \"494\u00d75=\"],\n  [\"918\u00d73=\", \"633\u00d74=\"],\n  [\"821\u00d75=\", \"553\u00d73=\"],\n  [\"906\u00d78=\", \"261\u00d73=\"],\n  [\"798\u00d74=\", \"850\u00d72=\"],\n  [\"843\u00d76=\", \"951\u00d78=\"],\n  [\"473\u00d72=\", \"476\u00d79=\"],\n  [\"606\u00d79=\", \"492\u00d79=\"],\n  [\"294\u00d78=\", \"738\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# worksheet table with its new problem, preserving formatting (the\n# Find/Replace operates on the text run in place, so font/size/\n# justification carried by the surrounding run/paragraph are untouched).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"520\u00d75=\", \"817\u00d78=\"),\n    @(\"249\u00d76=\", \"779\u00d79=\"),\n    @(\"442\u00d77=\", \"936\u00d78=\"),\n    @(\"377\u00d75=\", \"220\u00d76=\"),\n    @(\"382\u00d76=\", \"984\u00d76=\"),\n    @(\"468\u00d79=\", \"564\u00d72=\"),\n    @(\"610\u00d73=\", \"944\u00d72=\"),\n    @(\"374\u00d79=\", \"555\u00d75=\"),\n    @(\"171\u00d79=\", \"357\u00d76=\"),\n    @(\"966\u00d79=\", \"982\u00d74=\"),\n    @(\"942\u00d79=\", \"595\u00d76=\"),\n    @(\"710\u00d76=\", \"931\u00d75=\"),\n    @(\"135\u00d76=\", \"285\u00d77=\"),\n    @(\"576\u00d74=\", \"816\u00d73=\"),\n    @(\"429\u00d78=\", \"755\u00d75=\"),\n    @(\"383\u00d73=\", \"825\u00d74=\"),\n    @(\"502\u00d79=\", \"494\u00d75=\"),\n    @(\"918\u00d73=\", \"633\u00d74=\"),\n    @(\"821\u00d75=\", \"553\u00d73=\"),\n    @(\"906\u00d78=\", \"261\u00d73=\"),\n    @(\"798\u00d74=\", \"850\u00d72=\"),\n    @(\"843\u00d76=\", \"951\u00d78=\"),\n    @(\"473\u00d72=\", \"476\u00d79=\"),\n    @(\"606\u00d79=\", \"492\u00d79=\"),\n    @(\"294\u00d78=\", \"738\u00d79=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
